# LodeRunner FrameCompare - progress update on sheet "V3":
# fill in actual v3/v1 frame counts (columns B/C) for places 7-12 (rows 9-14),
# which previously had no data (column D formula was showing "-").
# Also move the active selection to C16 (next row to fill in).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V3")

$data = @(
    @(9,  23845, 28099),
    @(10, 27369, 32036),
    @(11, 32334, 37877),
    @(12, 35635, 41646),
    @(13, 38350, 45474),
    @(14, 41932, 50990)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

$ws.Activate()
$ws.Range("C16").Select()
